$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 319, shifting existing rows 319:470 down to 320:471
$ws.Rows.Item(319).Insert()

# Populate the new row 319 with the new weekly data point
$ws.Cells.Item(319, 1).Value = 9
$ws.Cells.Item(319, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(319, 3).Value = "Metropolitana"
$ws.Cells.Item(319, 4).Value = 45141
$ws.Cells.Item(319, 5).Value = 13
$ws.Cells.Item(319, 6).Value = 300000001
$ws.Cells.Item(319, 7).Value = "Rabanito"
$ws.Cells.Item(319, 8).Value = "Sin especificar"
$ws.Cells.Item(319, 9).Value = "Primera"
$ws.Cells.Item(319, 10).Value = 7000
$ws.Cells.Item(319, 11).Value = 3000
$ws.Cells.Item(319, 12).Value = 3000
$ws.Cells.Item(319, 13).Value = 3000
$ws.Cells.Item(319, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(319, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(319, 16).Value = 30
$ws.Cells.Item(319, 17).Value = 100
$ws.Cells.Item(319, 18).Value = "Hortaliza"

# Ensure date cell keeps same numeric format style as other D-column cells
$ws.Cells.Item(319, 4).NumberFormat = $ws.Cells.Item(320, 4).NumberFormat
